$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.687.76'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = '3.853.12'
$ws.Range("E3").Value = '  -1.46%  '
$ws.Range("E4").Value = '  +0.27%  '
$ws.Range("D5").Value = '''600.64'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.51%  '
$ws.Range("D6").Value = '''167.63'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.42%  '
$ws.Range("D7").Value = '3.856.57'
$ws.Range("E7").Value = '  -1.22%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '''0.529'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.66%  '
$ws.Range("D10").Value = '''0.165'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.91%  '
$ws.Range("D11").Value = '''6.32'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.30%  '
$ws.Range("D12").Value = '''0.462'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.37%  '
$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = '''37.48'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.98%  '
$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").Value = '''0.0000250'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.37%  '
$ws.Range("D15").Value = '4.492.54'
$ws.Range("E15").Value = '  -1.22%  '
$ws.Range("D16").Value = '3.827.15'
$ws.Range("E16").Value = '  -0.97%  '
$ws.Range("D17").Value = '68.787.95'
$ws.Range("E17").Value = '  -0.60%  '
$ws.Range("D18").Value = '''18.34'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.57%  '
$ws.Range("D19").Value = '''7.52'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.03%  '
$ws.Range("E20").Value = '  -1.21%  '
$ws.Range("D21").Value = '''10.80'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.81%  '
$ws.Range("D22").Value = '''478.29'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.32%  '
$ws.Range("D23").Value = '''0.732'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.00%  '
$ws.Range("D24").Value = '''0.0000161'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.77%  '
$ws.Range("D25").Value = '''84.62'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.12%  '
$ws.Range("D26").Value = '''2.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.55%  '
$ws.Range("D27").Value = '''12.22'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.36%  '
$ws.Range("D28").Value = '''10.01'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.22%  '
$ws.Range("E29").Value = '  -0.08%  '
$ws.Range("D30").Value = '''2.96'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("D31").Value = '3.997.19'
$ws.Range("E31").Value = '  -1.14%  '
$ws.Range("D32").Value = '''7.72'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.65%  '
$ws.Range("D33").Value = '''2.30'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.49%  '
$ws.Range("D34").Value = '''31.11'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.15%  '
$ws.Range("D35").Value = '3.815.58'
$ws.Range("E35").Value = '  -0.63%  '
$ws.Range("D36").Value = '''0.105'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.12%  '
$ws.Range("B37").Value = 'Filecoin'
$ws.Range("C37").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D37").Value = '''5.97'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.70%  '
$ws.Range("B38").Value = 'Kaspa'
$ws.Range("C38").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D38").Value = '''0.140'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.16%  '
$ws.Range("B39").Value = 'Mantle'
$ws.Range("C39").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D39").Value = '''1.01'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.43%  '
$ws.Range("B40").Value = 'dogwifhat'
$ws.Range("C40").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D40").Value = '''3.24'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +6.20%  '
$ws.Range("D41").Value = '''0.998'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.00%  '
$ws.Range("D42").Value = '''0.316'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.43%  '
$ws.Range("D43").Value = '''2.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.67%  '
$ws.Range("D44").Value = '''425.55'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -4.52%  '
$ws.Range("D45").Value = '''47.71'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.59%  '
$ws.Range("B46").Value = 'USDe'
$ws.Range("C46").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D46").Value = '''1.00'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("B47").Value = 'Cosmos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D47").Value = '''8.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.10%  '
$ws.Range("B48").Value = 'VeChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D48").Value = '''0.0358'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.08%  '
$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").Value = '''140.99'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.47%  '
$ws.Range("B50").Value = 'FLOKI'
$ws.Range("C50").Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range("D50").Value = '''0.000266'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +10.12%  '
$ws.Range("B51").Value = 'Arweave'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D51").Value = '''38.85'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.99%  '
